# Update the two-digit multiplication problems to the new set of values.
# Each cell contains a unique "AxB=" expression, so plain Find/Replace
# (not "Replace All") on each exact old string is safe and unambiguous.

$d = $word.ActiveDocument

$replacements = @(
    @("47×82=", "43×16="),
    @("12×70=", "17×24="),
    @("70×51=", "41×52="),
    @("33×26=", "26×68="),
    @("43×94=", "83×84="),
    @("19×30=", "89×47="),
    @("56×67=", "46×18="),
    @("83×61=", "69×98="),
    @("12×20=", "16×19="),
    @("20×70=", "35×69="),
    @("67×42=", "82×23="),
    @("45×77=", "58×89="),
    @("25×80=", "70×28="),
    @("76×34=", "67×55="),
    @("75×51=", "27×30="),
    @("26×93=", "94×36="),
    @("45×57=", "89×11="),
    @("58×13=", "88×73="),
    @("84×18=", "93×39="),
    @("63×28=", "73×23="),
    @("50×56=", "72×24="),
    @("43×34=", "81×79="),
    @("71×48=", "46×17="),
    @("35×76=", "34×28="),
    @("72×83=", "56×95=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Replaced $($replacements.Count) multiplication expressions"
